# Append 5 new data rows (rows 7-11) to the "Kedar Jadhav" sheet, duplicating
# existing rows 2, 4, 6, 3 and 5 (in that order), matching the committed diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# Each entry: venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr
$rows = @(
    @(" Dubai (DSC)", " October 02 2020", "Sunrisers won by 7 runs", "Chennai Super Kings", "Sunrisers Hyderabad", "Kedar Jadhav$nbsp", "3", "10", "0", "0", "30.00"),
    @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Chennai Super Kings", "Rajasthan Royals", "Kedar Jadhav$nbsp", "22", "16", "3", "0", "137.50"),
    @(" Abu Dhabi", " October 19 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Chennai Super Kings", "Rajasthan Royals", "Kedar Jadhav$nbsp", "4", "7", "0", "0", "57.14"),
    @(" Dubai (DSC)", " September 25 2020", "Capitals won by 44 runs", "Chennai Super Kings", "Delhi Capitals", "Kedar Jadhav$nbsp", "26", "21", "3", "0", "123.80"),
    @(" Abu Dhabi", " October 07 2020", "KKR won by 10 runs", "Chennai Super Kings", "Kolkata Knight Riders", "Kedar Jadhav$nbsp", "7", "12", "1", "0", "58.33")
)

$startRow = 7
$numberColumns = @("G", "H", "I", "J", "K")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]

    $ws.Range("A$r").Value = $values[0]
    $ws.Range("B$r").Value = $values[1]
    $ws.Range("C$r").Value = $values[2]
    $ws.Range("D$r").Value = $values[3]
    $ws.Range("E$r").Value = $values[4]
    $ws.Range("F$r").Value = $values[5]

    # Columns G-K look numeric but must stay stored as text, matching the
    # source workbook (t="str" cells with a numberStoredAsText ignored error).
    for ($c = 0; $c -lt $numberColumns.Count; $c++) {
        $col = $numberColumns[$c]
        $cell = $ws.Range("$col$r")
        $cell.NumberFormat = "@"
        $cell.Value = $values[6 + $c]
    }
}
